# Ajout de spécifications dans les scénarios
#
# In both "CRUD" scenario sheets, the "Ajouter" (Add) use-case row only
# required "Tout les champs sont remplis" (all fields are filled) as its
# condition. This commit tightens that condition to also require the
# entered data to be valid, matching the wording already used by the
# "Modifier" / "Supprimer" rows: "Tout les champs sont remplis et les
# données entrées sont valides".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("CRUD Locations")
$ws1.Range("B8").Value = "Tout les champs sont remplis et les données entrées sont valides"
$ws1.Rows.Item(8).RowHeight = 30

$ws2 = $wb.Worksheets.Item("CRUD Utilisateurs")
$ws2.Range("B8").Value = "Tout les champs sont remplis et les données entrées sont valides"
$ws2.Rows.Item(8).RowHeight = 30

# Leave the final selection/active sheet on "CRUD Utilisateurs" (matches
# the saved state in the updated workbook).
$ws1.Range("B8").Select()
$ws2.Activate()
$ws2.Range("C11").Select()
